# chore: adapt column header formatting to respective input file names (#7)
#
# Renames the "_old" / "_new" header-name suffixes to the concrete
# format-version identifiers "_FV2410" / "_FV2504", turns the sheet's
# data range into a proper Excel Table (ListObject) with an AutoFilter,
# and freezes the header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename the header row (row 1) --------------------------------
# Columns A:J carried the "..._old" suffix, column K is the untouched
# "diff" column, and columns L:U carried the "..._new" suffix.
$oldHeaders = @(
    "Segmentname_old",
    "Segmentgruppe_old",
    "Segment_old",
    "Datenelement_old",
    "Segment ID_old",
    "Code_old",
    "Qualifier_old",
    "Beschreibung_old",
    "Bedingungsausdruck_old",
    "Bedingung_old"
)

$newHeaders = @(
    "Segmentname_new",
    "Segmentgruppe_new",
    "Segment_new",
    "Datenelement_new",
    "Segment ID_new",
    "Code_new",
    "Qualifier_new",
    "Beschreibung_new",
    "Bedingungsausdruck_new",
    "Bedingung_new"
)

for ($i = 0; $i -lt $oldHeaders.Count; $i++) {
    $col = $i + 1
    $ws.Cells.Item(1, $col).Value = ($oldHeaders[$i] -replace "_old$", "_FV2410")
}

for ($i = 0; $i -lt $newHeaders.Count; $i++) {
    $col = $i + 12
    $ws.Cells.Item(1, $col).Value = ($newHeaders[$i] -replace "_new$", "_FV2504")
}

# --- 2. Turn the data range into an Excel Table ------------------------
$tableRange = $ws.Range("A1:U72")
$tbl = $ws.ListObjects.Add(1, $tableRange, $null, 1)
$tbl.Name = "Table1"

# --- 3. Freeze the header row ------------------------------------------
[void]$ws.Activate()
[void]$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
